# Update F-column values ("人气"/heat count) on the "展览" and "全部类型" sheets
# to reflect refreshed data, per commit: "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row -> new F value, for sheet "展览"
$zhanlanUpdates = @{
    4  = 8372
    5  = 8372
    9  = 7341
    11 = 511
    18 = 149
    19 = 12177
    22 = 2470
    23 = 3537
    27 = 111
    30 = 3353
    33 = 1719
    36 = 6018
    41 = 904
    49 = 20
    50 = 117
}

# Row -> new F value, for sheet "全部类型"
$quanbuUpdates = @{
    8  = 8372
    12 = 7341
    13 = 7341
    15 = 511
    21 = 149
    23 = 12177
    27 = 2470
    28 = 2470
    29 = 3537
    30 = 111
    34 = 3353
    37 = 1719
    40 = 6018
    46 = 904
    52 = 117
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
